$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": fill column E (rows 2-25) with actual numeric values ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$values = @{
    2 = 99.3
    3 = 87.44
    4 = 95.03
    5 = 72.94
    6 = 68.3
    7 = 61.53
    8 = 77
    9 = 98.69
    10 = 104.56
    11 = 60.67
    12 = 21.81
    13 = 1.17
    14 = 0
    15 = -0.01
    16 = 0
    17 = 2.71
    18 = 6.27
    19 = 52.91
    20 = 75.01000000000001
    21 = 114.64
    22 = 124.9
    23 = 111.6
    24 = 128.32
    25 = 99.01000000000001
}

foreach ($row in $values.Keys) {
    $ws1.Cells.Item($row, 5).Value = $values[$row]
}

# --- Sheet "CO2": insert a new row at position 2 for 2025-06-16 ---
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Rows.Item(2).Insert()
$ws3.Rows.Item(2).ClearFormats()

$ws3.Cells.Item(2, 1).NumberFormat = "@"
$ws3.Cells.Item(2, 1).Value = "2025-06-16"
$ws3.Cells.Item(2, 1).ClearFormats()

$ws3.Cells.Item(2, 2).Value = "-"
